$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("省份脚本")

$ws.Range("E3").Value = '上海市频道和卫视'
$ws.Range("E5").Value = '安徽省频道和卫视'
$ws.Range("E9").Value = '北京市频道和卫视'
$ws.Range("E11").Value = '福建省频道和卫视'
$ws.Range("E14").Value = '甘肃省频道和卫视'
$ws.Range("E16").Value = '广东省频道和卫视'
$ws.Range("E18").Value = '由腾讯云提供的广东省频道和卫视'
$ws.Range("E19").Value = '深圳市频道和卫视,部分广东频道'
$ws.Range("E21").Value = '贵州省频道和卫视'
$ws.Range("E22").Value = '由腾讯云提供的贵州省频道和卫视'
$ws.Range("E23").Value = '海南省频道和卫视'
$ws.Range("E24").Value = '河北省频道和卫视,河北地方频道'
$ws.Range("E26").Value = '河南省频道和卫视'
$ws.Range("E27").Value = '河南省频道和卫视,部分河南地方频道'
$ws.Range("E35").Value = '湖北省频道和卫视,央视和其他省卫视,劲爆体育'
$ws.Range("E36").Value = '湖北省频道和卫视'
$ws.Range("E39").Value = '湖南省频道和卫视'
$ws.Range("E41").Value = '吉林省频道和卫视'
$ws.Range("E44").Value = '江苏省频道和卫视,江苏地方频道'
$ws.Range("E46").Value = '由腾讯云提供的江苏省频道和卫视,江苏地方频道'
$ws.Range("E47").Value = '江西省频道和卫视'
$ws.Range("E52").Value = '辽宁省频道和卫视,辽宁地方频道'
$ws.Range("E54").Value = '内蒙古自治区和卫视,部分内蒙古地方频道'
$ws.Range("E56").Value = '由腾讯云提供的内蒙古自治区和卫视,部分内蒙古地方频道'
$ws.Range("E59").Value = '陕西省频道和卫视'
$ws.Range("E63").Value = '山东省频道和卫视'
$ws.Range("E64").Value = '山东省频道和卫视,部分山东地方频道'
$ws.Range("E68").Value = '山西省频道和卫视'
$ws.Range("E77").Value = '四川省频道和卫视'
$ws.Range("E80").Value = '西藏自治区频道和卫视'
$ws.Range("E82").Value = '新疆维吾尔自治区频道和卫视'
$ws.Range("E85").Value = '云南省频道和卫视'
$ws.Range("E87").Value = '浙江省频道和卫视'
$ws.Range("E92").Value = '香港卫视'
